$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "So, the output from the last FC layer " -> bold "the last FC layer"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("So, the output from the last FC layer ") | Out-Null
$start1 = $r1.Start
$end1 = $r1.End
$bold1 = $d.Range($start1 + 20, $end1 - 1)
$bold1.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: "and the output of the last pooling layer is passed to " ->
#           bold "the last pooling layer"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("and the output of the last pooling layer is passed to ") | Out-Null
$start2 = $r2.Start
$bold2 = $d.Range($start2 + 18, $start2 + 18 + 22)
$bold2.Bold = 1

# ---------------------------------------------------------------------------
# Change 3: restructure "The input to a regressor is a set of N pairs (P..."
#           into three paragraphs.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("The input to a regressor is a set of N pairs (P") | Out-Null
$insertPoint3 = $d.Range($r3.Start, $r3.Start)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">We have class specific </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>regressors</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">. Once the SVMs give the output class for a proposed region, we use the class-specific </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>regressor</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> to predict the bounding box.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/></w:pPr></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint3.InsertXML($xml3)

$r3b = $d.Content
$r3b.Find.Execute("The input to a regressor is a set of N pairs (P") | Out-Null
$r3b.Text = "The input is a set of N pairs (P"

# ---------------------------------------------------------------------------
# Change 4: add <w:lastRenderedPageBreak/> before "Similarly, during testing",
#           remove the manual page break run, and add three empty paragraphs.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Similarly, during testing") | Out-Null
$insertPoint4 = $d.Range($r4.Start, $r4.Start)
$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint4.InsertXML($xml4)

$r4b = $d.Content
$r4b.Find.Execute("output, which is relative to ground-truth box.") | Out-Null
$breakChar = $d.Range($r4b.End, $r4b.End + 1)
$breakChar.Delete()

$para4 = $r4b.Paragraphs.First
$endOfPara4 = $para4.Range.End
$insertPoint4b = $d.Range($endOfPara4, $endOfPara4)
$xml4b = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p/><w:p/><w:p/>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint4b.InsertXML($xml4b)

Write-Output "All changes applied"
